$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.228626251220703
$ws.Range("B1").Value = 2.447827577590942
$ws.Range("C1").Value = 2.566384077072144
$ws.Range("D1").Value = 3.405359983444214
$ws.Range("E1").Value = 1.641739010810852
